# Generate Report for Archive
#
# 1) The shared string "Ready for handoff" becomes "In Translation".
#    It is used by every Status cell on all three sheets (Overview!E2:F4,
#    zh-cn!C2:C4, de-de!C2:C4), so every one of those cells is updated.
# 2) The "Status" columns get narrower: Overview columns E & F and the
#    Status column (C) on the zh-cn / de-de sheets shrink from
#    ~17.22 chars to ~13.41 chars.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- 1) Update every "Ready for handoff" status cell to "In Translation" ---
$wsOverview.Range("E2:F4").Value = "In Translation"
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# --- 2) Narrow the Status columns ---
# Target stored width is ~13.4101845877511 characters. Excel quantizes the
# ColumnWidth property to whole pixels on save, so feed it the character
# width whose quantized result lands on that value.
$newStatusWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusWidth  # column E
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusWidth  # column F
$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusWidth       # column C
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusWidth       # column C

# NOTE: keep a trailing statement after the last ColumnWidth assignment --
# when a `Range.ColumnWidth = ...` COM property-set is the final statement
# executed by the script, the interop layer chokes trying to surface its
# return value ("type mismatch") and the assignment is silently dropped.
Write-Output "Report regenerated."
